$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.048.42'
$ws.Range("E2").Value = '  +2.17%  '
$ws.Range("D3").Value = '1.673.30'
$ws.Range("E3").Value = '  +2.80%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.21'
$ws.Range("E5").Value = '  +1.33%  '
$ws.Range("E6").Value = '  +1.69%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  +2.07%  '
$ws.Range("E9").Value = '  +1.11%  '
$ws.Range("E10").Value = '  +4.86%  '
$ws.Range("E11").Value = '  +4.75%  '
$ws.Range("D12").Value = '1.909.32'
$ws.Range("E12").Value = '  +3.05%  '
$ws.Range("D13").Value = '1.674.08'
$ws.Range("E13").Value = '  +2.69%  '
$ws.Range("E14").Value = '  +1.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '65.80'
$ws.Range("E15").Value = '  +2.84%  '
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("D17").Value = '27.059.73'
$ws.Range("E17").Value = '  +2.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '235.09'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("E19").Value = '  +1.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.72'
$ws.Range("E20").Value = '  -1.59%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  +3.16%  '
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("E24").Value = '  +1.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.11'
$ws.Range("E25").Value = '  -1.45%  '
$ws.Range("E26").Value = '  +0.99%  '
$ws.Range("E27").Value = '  +0.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.92'
$ws.Range("E28").Value = '  +1.92%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("E31").Value = '  +1.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.33'
$ws.Range("E32").Value = '  +2.13%  '
$ws.Range("D33").Value = '1.458.21'
$ws.Range("E33").Value = '  -4.38%  '
$ws.Range("E34").Value = '  +5.36%  '
$ws.Range("E35").Value = '  +5.52%  '
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.572'
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.894'
$ws.Range("E38").Value = '  +6.84%  '
$ws.Range("E39").Value = '  +1.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.07'
$ws.Range("E40").Value = '  +3.21%  '
$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.02'
$ws.Range("E41").Value = '  +11.49%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  +3.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.01'
$ws.Range("E44").Value = '  +4.83%  '
$ws.Range("D45").Value = '1.817.89'
$ws.Range("E45").Value = '  +3.04%  '
$ws.Range("E46").Value = '  +1.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.31'
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("E48").Value = '  +1.51%  '
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.100'
$ws.Range("E50").Value = '  +3.82%  '
$ws.Range("E51").Value = '  +1.34%  '
